$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "36÷5="
$t.Cell(1, 2).Range.Text = "13÷7="
$t.Cell(1, 3).Range.Text = "73÷3="
$t.Cell(1, 4).Range.Text = "76÷2="
$t.Cell(1, 5).Range.Text = "64÷9="
$t.Cell(5, 1).Range.Text = "52÷3="
$t.Cell(5, 2).Range.Text = "76÷3="
$t.Cell(5, 3).Range.Text = "57÷7="
$t.Cell(5, 4).Range.Text = "67÷5="
$t.Cell(5, 5).Range.Text = "26÷7="
$t.Cell(9, 1).Range.Text = "52÷2="
$t.Cell(9, 2).Range.Text = "32÷8="
$t.Cell(9, 3).Range.Text = "80÷4="
$t.Cell(9, 4).Range.Text = "19÷5="
$t.Cell(9, 5).Range.Text = "90÷7="
$t.Cell(13, 1).Range.Text = "32÷2="
$t.Cell(13, 2).Range.Text = "12÷3="
$t.Cell(13, 3).Range.Text = "99÷9="
$t.Cell(13, 4).Range.Text = "91÷2="
$t.Cell(13, 5).Range.Text = "97÷9="
$t.Cell(17, 1).Range.Text = "57÷9="
$t.Cell(17, 2).Range.Text = "70÷6="
$t.Cell(17, 3).Range.Text = "22÷5="
$t.Cell(17, 4).Range.Text = "86÷2="
$t.Cell(17, 5).Range.Text = "45÷6="
